# Auto-generated Excel COM-interop edit script
# Applies the cryptos.xlsx data refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.952.48"
$ws.Range("E2").Value = "  +0.22%  "

# Row 3
$ws.Range("D3").Value = "2.751.71"
$ws.Range("E3").Value = "  -0.33%  "

# Row 4
$ws.Range("E4").Value = "  +0.18%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.33%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.99%  "

# Row 7
$ws.Range("E7").Value = "  +0.35%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.605"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.96%  "

# Row 9
$ws.Range("E9").Value = "  -2.13%  "

# Row 10
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.384"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.49%  "

# Row 11
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.159"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.88%  "

# Row 12
$ws.Range("E12").Value = "  -15.80%  "

# Row 13
$ws.Range("D13").Value = "3.237.98"
$ws.Range("E13").Value = "  -0.08%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.16%  "

# Row 15
$ws.Range("D15").Value = "63.691.93"
$ws.Range("E15").Value = "  -0.05%  "

# Row 16
$ws.Range("E16").Value = "  -1.77%  "

# Row 17
$ws.Range("D17").Value = "2.758.11"
$ws.Range("E17").Value = "  +0.17%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.21%  "

# Row 19
$ws.Range("E19").Value = "  -0.43%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "357.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.86%  "

# Row 21
$ws.Range("E21").Value = "  -2.67%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.543"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.58%  "

# Row 23
$ws.Range("E23").Value = "  -0.13%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.83%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.171"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.76%  "

# Row 26
$ws.Range("E26").Value = "  +1.06%  "

# Row 27
$ws.Range("E27").Value = "  -1.18%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0903"
$ws.Range("E28").Value = "  -1.33%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.28%  "

# Row 30
$ws.Range("E30").Value = "  -2.37%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "170.87"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.40%  "

# Row 32
$ws.Range("E32").Value = "  -2.03%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.27"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.14%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.58%  "

# Row 36
$ws.Range("E36").Value = "  +1.07%  "

# Row 37
$ws.Range("E37").Value = "  -1.15%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.982"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.53%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.97%  "

# Row 40
$ws.Range("E40").Value = "  -2.01%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "326.09"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.46%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.52%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.51%  "

# Row 44
$ws.Range("E44").Value = "  -0.11%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.84%  "

# Row 46
$ws.Range("E46").Value = "  -0.44%  "

# Row 47
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.630"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.28%  "

# Row 48
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.23%  "

# Row 49
$ws.Range("E49").Value = "  -0.10%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.29%  "

# Row 51
$ws.Range("E51").Value = "  +0.63%  "

